# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update case counters for several countries (Excel COM model keeps the
#   table sorted by "Casos totales" descending, so a couple of countries
#   change rank/row as their totals are refreshed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 17:35"

# 2) Estados Unidos (row 4) - counts refreshed, stays #1
$ws.Range("B4").Value = 1648961
$ws.Range("C4").Value = 3867
$ws.Range("D4").Value = 403315
$ws.Range("E4").Value = 1147871
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 128
$ws.Range("H4").Value = 97775

# 3) Row 29 - counts refreshed in place (no rank change)
$ws.Range("D29").Value = 13882
$ws.Range("E29").Value = 17163

# 4) Row 50 - counts refreshed in place (no rank change)
$ws.Range("D50").Value = 3530
$ws.Range("E50").Value = 6680
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 439

# 5) Argelia / Kazajistan swap ranks (rows 56-57).
#    Kazajistan's own numbers are unchanged, it just drops to row 57.
#    Argelia's numbers are refreshed and it rises to row 56.
$ws.Range("A56").Value = "Argelia"
$ws.Range("B56").Value = 8113
$ws.Range("C56").Value = 195
$ws.Range("D56").Value = 4426
$ws.Range("E56").Value = 3095
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 10
$ws.Range("H56").Value = 592

$ws.Range("A57").Value = "Kazajistan"
$ws.Range("B57").Value = 7919
$ws.Range("C57").Value = 322
$ws.Range("D57").Value = 4096
$ws.Range("E57").Value = 3788
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 35

# 6) Row 63 - counts refreshed in place (no rank change)
$ws.Range("B63").Value = 6994
$ws.Range("C63").Value = 147
$ws.Range("E63").Value = 3300
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 242

# 7) Row 79 - counts refreshed in place (no rank change)
$ws.Range("B79").Value = 2876
$ws.Range("C79").Value = 2
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 171

# 8) "Republica de Africa Central" rises to row 131 (with refreshed
#    numbers), pushing Jamaica / Tanzania / Etiopia / Madagascar down one
#    row each (their own numbers are unchanged). Congo (row 136) keeps its
#    row and numbers untouched.
$ws.Range("A131").Value = "Republica de Africa Central"
$ws.Range("B131").Value = 552
$ws.Range("C131").Value = 73
$ws.Range("D131").Value = 18
$ws.Range("E131").Value = 533
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 1

$ws.Range("A132").Value = "Jamaica"
$ws.Range("B132").Value = 544
$ws.Range("C132").Value = 10
$ws.Range("D132").Value = 191
$ws.Range("E132").Value = 344
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 9

$ws.Range("A133").Value = "Tanzania"
$ws.Range("B133").Value = 509
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 183
$ws.Range("E133").Value = 305
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 21

$ws.Range("A134").Value = "Etiopia"
$ws.Range("B134").Value = 494
$ws.Range("C134").Value = 61
$ws.Range("D134").Value = 151
$ws.Range("E134").Value = 338
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 5

$ws.Range("A135").Value = "Madagascar"
$ws.Range("B135").Value = 488
$ws.Range("C135").Value = 40
$ws.Range("D135").Value = 138
$ws.Range("E135").Value = 348
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 2
